$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Cells.Item(18, 8).Value = 111499.11
$ws.Cells.Item(18, 9).Value = 111499.11
$ws.Cells.Item(18, 11).Value = 111499.11
$ws.Cells.Item(18, 13).Value = -111215.11

# Row 19
$ws.Cells.Item(19, 8).Value = 10206226
$ws.Cells.Item(19, 9).Value = 7825968
$ws.Cells.Item(19, 10).Value = 14286668
$ws.Cells.Item(19, 11).Value = 7825968
$ws.Cells.Item(19, 12).Value = 14286668
$ws.Cells.Item(19, 13).Value = -7825793
$ws.Cells.Item(19, 14).Value = -14287018

# Row 76
$ws.Cells.Item(76, 8).Value = 3881.7856
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 3881.7856
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).ClearContents()
$ws.Cells.Item(76, 13).Value = 3881.7856
$ws.Cells.Item(76, 14).Value = -4511.7856

# Row 79
$ws.Cells.Item(79, 8).Value = 3881.7856
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 3881.7856
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).ClearContents()
$ws.Cells.Item(79, 13).Value = 3881.7856
$ws.Cells.Item(79, 14).Value = -6065.7856

# Row 80
$ws.Cells.Item(80, 8).Value = 299.5
$ws.Cells.Item(80, 9).Value = 316
$ws.Cells.Item(80, 10).Value = 288.5
$ws.Cells.Item(80, 11).Value = 948
$ws.Cells.Item(80, 12).Value = 865.5
$ws.Cells.Item(80, 13).Value = 50
$ws.Cells.Item(80, 14).Value = -2861.5

# Row 83
$ws.Cells.Item(83, 8).Value = 299.5
$ws.Cells.Item(83, 9).Value = 316
$ws.Cells.Item(83, 10).Value = 288.5
$ws.Cells.Item(83, 11).Value = 2844
$ws.Cells.Item(83, 12).Value = 2596.5
$ws.Cells.Item(83, 13).Value = 2148
$ws.Cells.Item(83, 14).Value = -12580.5

# Row 125
$ws.Cells.Item(125, 8).Value = 5175.125
$ws.Cells.Item(125, 9).Value = 5000
$ws.Cells.Item(125, 10).Value = 5200.143
$ws.Cells.Item(125, 11).Value = 45000
$ws.Cells.Item(125, 12).Value = 46801.287
$ws.Cells.Item(125, 13).Value = -42540
$ws.Cells.Item(125, 14).Value = -51721.287

# Row 138
$ws.Cells.Item(138, 8).Value = 2758.4644
$ws.Cells.Item(138, 9).Value = 1194.8148
$ws.Cells.Item(138, 10).Value = 4214.276
$ws.Cells.Item(138, 11).Value = 3584.4444
$ws.Cells.Item(138, 12).Value = 12642.828
$ws.Cells.Item(138, 13).Value = 1555.5556
$ws.Cells.Item(138, 14).Value = -22922.828

# Row 140
$ws.Cells.Item(140, 8).Value = 52429.5
$ws.Cells.Item(140, 10).Value = 52429.5
$ws.Cells.Item(140, 12).Value = 52429.5
$ws.Cells.Item(140, 14).Value = -62789.5

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Cells.Item(4, 8).Value = 454.2857
$ws.Cells.Item(4, 9).Value = 454.2857
$ws.Cells.Item(4, 11).Value = 454.2857
$ws.Cells.Item(4, 13).Value = -338.2857

# Row 32
$ws.Cells.Item(32, 8).Value = 10003570
$ws.Cells.Item(32, 9).Value = 13159309
$ws.Cells.Item(32, 10).Value = 10398.042
$ws.Cells.Item(32, 11).Value = 13159309
$ws.Cells.Item(32, 12).Value = 10398.042
$ws.Cells.Item(32, 13).Value = -13159022
$ws.Cells.Item(32, 14).Value = -10972.042

# Row 45
$ws.Cells.Item(45, 8).Value = 1781.1613
$ws.Cells.Item(45, 9).Value = 1780.0667
$ws.Cells.Item(45, 11).Value = 1780.0667
$ws.Cells.Item(45, 13).Value = -1403.0667

# Row 122
$ws.Cells.Item(122, 8).Value = 964.875
$ws.Cells.Item(122, 9).Value = 948.6667
$ws.Cells.Item(122, 11).Value = 2846.0001
$ws.Cells.Item(122, 13).Value = -396.0001000000002

# Row 132
$ws.Cells.Item(132, 8).Value = 2805.5833
$ws.Cells.Item(132, 9).Value = 2432.353
$ws.Cells.Item(132, 10).Value = 3712
$ws.Cells.Item(132, 11).Value = 7297.059
$ws.Cells.Item(132, 12).Value = 11136
$ws.Cells.Item(132, 13).Value = -4767.059
$ws.Cells.Item(132, 14).Value = -16196

$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 13).ClearContents()

# Row 19
$ws.Cells.Item(19, 8).Value = 19000
$ws.Cells.Item(19, 10).Value = 19000
$ws.Cells.Item(19, 12).Value = 19000
$ws.Cells.Item(19, 14).Value = -19346

# Row 86
$ws.Cells.Item(86, 8).Value = 2331.5
$ws.Cells.Item(86, 9).Value = 1754.2727
$ws.Cells.Item(86, 10).Value = 3601.4
$ws.Cells.Item(86, 11).Value = 1754.2727
$ws.Cells.Item(86, 12).Value = 3601.4
$ws.Cells.Item(86, 13).Value = -631.2727
$ws.Cells.Item(86, 14).Value = -5847.4

# Row 89
$ws.Cells.Item(89, 8).Value = 2331.5
$ws.Cells.Item(89, 9).Value = 1754.2727
$ws.Cells.Item(89, 10).Value = 3601.4
$ws.Cells.Item(89, 11).Value = 8771.363499999999
$ws.Cells.Item(89, 12).Value = 18007
$ws.Cells.Item(89, 13).Value = -3155.363499999999
$ws.Cells.Item(89, 14).Value = -29239

# Row 105
$ws.Cells.Item(105, 8).Value = 2832.16
$ws.Cells.Item(105, 9).Value = 1600
$ws.Cells.Item(105, 10).Value = 2969.0667
$ws.Cells.Item(105, 11).Value = 1600
$ws.Cells.Item(105, 12).Value = 2969.0667
$ws.Cells.Item(105, 13).Value = 147
$ws.Cells.Item(105, 14).Value = -6463.066699999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 330.27274
$ws.Cells.Item(22, 9).Value = 321
$ws.Cells.Item(22, 10).Value = 372
$ws.Cells.Item(22, 11).Value = 321
$ws.Cells.Item(22, 12).Value = 372
$ws.Cells.Item(22, 13).Value = 29
$ws.Cells.Item(22, 14).Value = -1072

# Row 62
$ws.Cells.Item(62, 8).Value = 168335500
$ws.Cells.Item(62, 9).Value = 168335500
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 168335500
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).Value = -168334876

# Row 65
$ws.Cells.Item(65, 8).Value = 168335500
$ws.Cells.Item(65, 9).Value = 168335500
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 841677500
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).Value = -841674380

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 1011.4828
$ws.Cells.Item(131, 10).Value = 1131.8695
$ws.Cells.Item(131, 12).Value = 3395.6085
$ws.Cells.Item(131, 14).Value = -13475.6085

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Cells.Item(5, 8).Value = 10250
$ws.Cells.Item(5, 10).Value = 12000
$ws.Cells.Item(5, 12).Value = 12000
$ws.Cells.Item(5, 14).Value = -12224

# Row 18
$ws.Cells.Item(18, 8).Value = 3335666.8
$ws.Cells.Item(18, 9).Value = 3000
$ws.Cells.Item(18, 10).Value = 5002000
$ws.Cells.Item(18, 11).Value = 3000
$ws.Cells.Item(18, 12).Value = 5002000
$ws.Cells.Item(18, 13).Value = -2707
$ws.Cells.Item(18, 14).Value = -5002586

# Row 21
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 13).ClearContents()

# Row 30
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 13).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 3664.889
$ws.Cells.Item(7, 9).Value = 3426.2856
$ws.Cells.Item(7, 10).Value = 4500
$ws.Cells.Item(7, 11).Value = 3426.2856
$ws.Cells.Item(7, 12).Value = 4500
$ws.Cells.Item(7, 13).Value = -3314.2856
$ws.Cells.Item(7, 14).Value = -4724

# Row 20
$ws.Cells.Item(20, 8).Value = 15750
$ws.Cells.Item(20, 10).Value = 15750
$ws.Cells.Item(20, 12).Value = 15750
$ws.Cells.Item(20, 14).Value = -16202

# Row 21
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 13).ClearContents()

# Row 68
$ws.Cells.Item(68, 8).Value = 177241.28
$ws.Cells.Item(68, 9).Value = 528078.2
$ws.Cells.Item(68, 10).Value = 1822.8158
$ws.Cells.Item(68, 11).Value = 528078.2
$ws.Cells.Item(68, 12).Value = 1822.8158
$ws.Cells.Item(68, 13).Value = -527329.2
$ws.Cells.Item(68, 14).Value = -3320.8158

# Row 71
$ws.Cells.Item(71, 8).Value = 177241.28
$ws.Cells.Item(71, 9).Value = 528078.2
$ws.Cells.Item(71, 10).Value = 1822.8158
$ws.Cells.Item(71, 11).Value = 2640391
$ws.Cells.Item(71, 12).Value = 9114.079
$ws.Cells.Item(71, 13).Value = -2636647
$ws.Cells.Item(71, 14).Value = -16602.079

# Row 122
$ws.Cells.Item(122, 8).Value = 3478.5454
$ws.Cells.Item(122, 9).Value = 3152
$ws.Cells.Item(122, 10).Value = 3551.111
$ws.Cells.Item(122, 11).Value = 9456
$ws.Cells.Item(122, 12).Value = 10653.333
$ws.Cells.Item(122, 13).Value = -7006
$ws.Cells.Item(122, 14).Value = -15553.333

# Row 126
$ws.Cells.Item(126, 8).Value = 3664.889
$ws.Cells.Item(126, 9).Value = 3426.2856
$ws.Cells.Item(126, 10).Value = 4500
$ws.Cells.Item(126, 11).Value = 10278.8568
$ws.Cells.Item(126, 12).Value = 13500
$ws.Cells.Item(126, 13).Value = -7808.856800000001
$ws.Cells.Item(126, 14).Value = -18440

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Cells.Item(15, 8).Value = 10789.474
$ws.Cells.Item(15, 10).Value = 13750
$ws.Cells.Item(15, 12).Value = 13750
$ws.Cells.Item(15, 14).Value = -14326

# Row 107
$ws.Cells.Item(107, 8).Value = 412.7143
$ws.Cells.Item(107, 9).Value = 429.84616
$ws.Cells.Item(107, 10).Value = 190
$ws.Cells.Item(107, 11).Value = 1289.53848
$ws.Cells.Item(107, 12).Value = 570
$ws.Cells.Item(107, 13).Value = 630.4615200000001
$ws.Cells.Item(107, 14).Value = -4410
